# Updates the crypto price/volume snapshot on Sheet1 (columns D "Price" and
# E "Volume(1h)") for rows 2-51 to the latest scraped values, matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# All D/E cells in this sheet are stored as *text* (prices use a dotted
# European-style grouping, e.g. "62.043.49", and percentages keep their
# original padding/sign formatting, e.g. "  -0.36%  "). Writing a plain
# numeric-looking string straight into Range.Value lets Excel's COM layer
# auto-coerce it into a Double, which would silently reformat values such
# as "5.20" -> 5.2 or "152.40" -> 152.4. To avoid that, cells whose new
# text parses as a plain number are pre-formatted as Text ("@") before the
# value is written; everything else (percentages, and the few multi-dot
# prices Excel never parses as numeric anyway) is written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text looks like a plain number and must be forced to Text
# format first so Excel does not coerce them into a Double (which would
# normalize formatting, e.g. drop a trailing zero: "5.20" -> 5.2).
$textForced = @{
    "D5" = "561.78"
    "D6" = "143.53"
    "D12" = "5.20"
    "D14" = "26.13"
    "D19" = "11.25"
    "D20" = "322.67"
    "D23" = "0.999"
    "D24" = "67.34"
    "D27" = "560.21"
    "D31" = "8.19"
    "D39" = "5.46"
    "D40" = "152.40"
    "D44" = "2.25"
    "D45" = "147.32"
    "D46" = "3.63"
    "D48" = "19.90"
    "D49" = "0.595"
    "D50" = "0.0922"
    "D51" = "0.0229"
}
foreach ($addr in $textForced.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textForced[$addr]
}

# Remaining cells: either plain percentage text (never number-like) or
# price text containing 2+ dots (e.g. "62.043.49"), which Excel already
# keeps as text without any extra coaxing.
$plainValues = @{
    "D2" = "62.043.49"
    "E2" = "  -0.36%  "
    "D3" = "2.420.66"
    "E3" = "  -0.05%  "
    "E4" = "  +0.03%  "
    "E5" = "  -0.08%  "
    "E6" = "  -0.35%  "
    "E7" = "  -0.03%  "
    "E8" = "  -0.36%  "
    "D9" = "2.419.19"
    "E9" = "  +0.01%  "
    "E10" = "  -0.40%  "
    "E11" = "  +0.24%  "
    "E12" = "  -3.15%  "
    "E13" = "  -1.26%  "
    "E14" = "  +0.88%  "
    "E15" = "  -2.05%  "
    "D16" = "2.859.70"
    "E16" = "  -0.02%  "
    "D17" = "62.044.90"
    "E17" = "  -0.19%  "
    "D18" = "2.410.85"
    "E18" = "  -0.44%  "
    "E19" = "  -0.58%  "
    "E20" = "  -0.25%  "
    "E21" = "  -1.65%  "
    "E22" = "  +0.99%  "
    "E23" = "  -0.16%  "
    "E25" = "  +1.85%  "
    "E26" = "  -2.78%  "
    "E27" = "  -3.85%  "
    "D28" = "2.539.05"
    "E28" = "  -0.12%  "
    "E29" = "  +0.01%  "
    "D30" = "0.0₃0931"
    "E30" = "  -0.96%  "
    "E31" = "  -0.60%  "
    "E32" = "  -4.62%  "
    "E33" = "  -2.26%  "
    "E34" = "  -1.17%  "
    "E35" = "  -2.64%  "
    "E37" = "  -1.00%  "
    "E38" = "  -1.19%  "
    "E39" = "  -4.24%  "
    "E40" = "  -0.03%  "
    "E41" = "  -0.01%  "
    "E42" = "  -1.13%  "
    "E43" = "  +0.30%  "
    "E44" = "  -2.83%  "
    "E45" = "  -1.85%  "
    "E46" = "  -0.49%  "
    "E47" = "  -1.80%  "
    "E48" = "  -2.02%  "
    "E49" = "  +0.18%  "
    "E50" = "  -0.33%  "
    "E51" = "  -0.07%  "
}
foreach ($addr in $plainValues.Keys) {
    $ws.Range($addr).Value = $plainValues[$addr]
}

